# AFDP-7308 Combine Transcribe and OCR processing into a single media
# processing module.
#
# Renames the OCR-specific business process model / workflow references
# in the "OCR Workflow Rules" rule table (Sheet1) to the new, combined
# "Media Engine" naming, and updates the sheet's saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# $model: OCRBusinessProcessModel  ->  $model: MediaEngineBusinessProcessModel
$ws.Range("C14").Value = "`$model: MediaEngineBusinessProcessModel"

# com.armedia.acm.ocr.model.OCRBusinessProcessModel
#   -> com.armedia.acm.services.mediaengine.model.MediaEngineBusinessProcessModel
$ws.Range("D3").Value = "com.armedia.acm.services.mediaengine.model.MediaEngineBusinessProcessModel"

# OCRWorkflow -> MediaEngineWorkFlow
$ws.Range("E17").Value = "MediaEngineWorkFlow"
$ws.Range("E18").Value = "MediaEngineWorkFlow"

# Update the saved scroll position / selection on the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E20").Select()
